$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new block as Text so the "1"/"2"/"3" survey codes are stored
# as shared strings (matching the existing response rows) rather than numbers.
$ws.Range("A25:BW27").NumberFormat = "@"

# Row 25
$ws.Range("A25").Value = "3"
$ws.Range("B25").Value = "1"
$ws.Range("C25").Value = "1"
$ws.Range("D25").Value = "3"
$ws.Range("E25").Value = "3"
$ws.Range("F25").Value = "3"
$ws.Range("G25").Value = "1"
$ws.Range("H25").Value = "2"
$ws.Range("I25").Value = "3"
$ws.Range("J25").Value = "3"
$ws.Range("K25").Value = "1"
$ws.Range("L25").Value = "1"
$ws.Range("M25").Value = "2"
$ws.Range("N25").Value = "2"
$ws.Range("O25").Value = "2"
$ws.Range("P25").Value = "2"
$ws.Range("Q25").Value = "1"
$ws.Range("R25").Value = "2"
$ws.Range("S25").Value = "2"
$ws.Range("T25").Value = "2"
$ws.Range("U25").Value = "2"
$ws.Range("V25").Value = "1"
$ws.Range("W25").Value = "2"
$ws.Range("X25").Value = "2"
$ws.Range("Y25").Value = "2"
$ws.Range("Z25").Value = "3"
$ws.Range("AA25").Value = "1"
$ws.Range("AB25").Value = "3"
$ws.Range("AC25").Value = "3"
$ws.Range("AD25").Value = "2"
$ws.Range("AE25").Value = "3"
$ws.Range("AF25").Value = "1"
$ws.Range("AG25").Value = "2"
$ws.Range("AH25").Value = "3"
$ws.Range("AI25").Value = "3"
$ws.Range("AJ25").Value = "3"
$ws.Range("AK25").Value = "1"
$ws.Range("AL25").Value = "3"
$ws.Range("AM25").Value = "2"
$ws.Range("AN25").Value = "3"
$ws.Range("AO25").Value = "1"
$ws.Range("AP25").Value = "1"
$ws.Range("AQ25").Value = "1"
$ws.Range("AR25").Value = "1"
$ws.Range("AS25").Value = "1"
$ws.Range("AT25").Value = "3"
$ws.Range("AU25").Value = "1"
$ws.Range("AV25").Value = "3"
$ws.Range("AW25").Value = "3"
$ws.Range("AX25").Value = "3"
$ws.Range("AY25").Value = "2"
$ws.Range("AZ25").Value = "1"
$ws.Range("BA25").Value = "2"
$ws.Range("BB25").Value = "2"
$ws.Range("BC25").Value = "2"
$ws.Range("BD25").Value = "3"
$ws.Range("BE25").Value = "3"
$ws.Range("BF25").Value = "3"
$ws.Range("BG25").Value = "3"
$ws.Range("BH25").Value = "3"
$ws.Range("BI25").Value = "1"
$ws.Range("BJ25").Value = "1"
$ws.Range("BK25").Value = "1"
$ws.Range("BL25").Value = "1"
$ws.Range("BM25").Value = "1"
$ws.Range("BN25").Value = "1"
$ws.Range("BO25").Value = "1"
$ws.Range("BP25").Value = "3"
$ws.Range("BQ25").Value = "2"
$ws.Range("BR25").Value = "3"
$ws.Range("BS25").Value = "3"
$ws.Range("BT25").Value = "3"
$ws.Range("BU25").Value = "2"
$ws.Range("BV25").Value = "3"
$ws.Range("BW25").Value = "3"

# Row 26
$ws.Range("A26").Value = "1"
$ws.Range("B26").Value = "1"
$ws.Range("C26").Value = "3"
$ws.Range("D26").Value = "3"
$ws.Range("E26").Value = "3"
$ws.Range("F26").Value = "1"
$ws.Range("G26").Value = "1"
$ws.Range("H26").Value = "2"
$ws.Range("I26").Value = "1"
$ws.Range("J26").Value = "2"
$ws.Range("K26").Value = "2"
$ws.Range("L26").Value = "1"
$ws.Range("M26").Value = "2"
$ws.Range("N26").Value = "2"
$ws.Range("O26").Value = "2"
$ws.Range("P26").Value = "2"
$ws.Range("Q26").Value = "1"
$ws.Range("R26").Value = "2"
$ws.Range("S26").Value = "2"
$ws.Range("T26").Value = "2"
$ws.Range("U26").Value = "2"
$ws.Range("V26").Value = "1"
$ws.Range("W26").Value = "2"
$ws.Range("X26").Value = "2"
$ws.Range("Y26").Value = "2"
$ws.Range("Z26").Value = "1"
$ws.Range("AA26").Value = "1"
$ws.Range("AB26").Value = "3"
$ws.Range("AC26").Value = "1"
$ws.Range("AD26").Value = "3"
$ws.Range("AE26").Value = "2"
$ws.Range("AF26").Value = "1"
$ws.Range("AG26").Value = "2"
$ws.Range("AH26").Value = "2"
$ws.Range("AI26").Value = "2"
$ws.Range("AJ26").Value = "1"
$ws.Range("AK26").Value = "1"
$ws.Range("AL26").Value = "3"
$ws.Range("AM26").Value = "1"
$ws.Range("AN26").Value = "3"
$ws.Range("AO26").Value = "2"
$ws.Range("AP26").Value = "1"
$ws.Range("AQ26").Value = "2"
$ws.Range("AR26").Value = "2"
$ws.Range("AS26").Value = "2"
$ws.Range("AT26").Value = "3"
$ws.Range("AU26").Value = "1"
$ws.Range("AV26").Value = "3"
$ws.Range("AW26").Value = "3"
$ws.Range("AX26").Value = "3"
$ws.Range("AY26").Value = "3"
$ws.Range("AZ26").Value = "1"
$ws.Range("BA26").Value = "2"
$ws.Range("BB26").Value = "3"
$ws.Range("BC26").Value = "3"
$ws.Range("BD26").Value = "2"
$ws.Range("BE26").Value = "1"
$ws.Range("BF26").Value = "2"
$ws.Range("BG26").Value = "2"
$ws.Range("BH26").Value = "2"
$ws.Range("BI26").Value = "1"
$ws.Range("BJ26").Value = "1"
$ws.Range("BK26").Value = "3"
$ws.Range("BL26").Value = "1"
$ws.Range("BM26").Value = "3"
$ws.Range("BN26").Value = "1"
$ws.Range("BO26").Value = "1"
$ws.Range("BP26").Value = "3"
$ws.Range("BQ26").Value = "3"
$ws.Range("BR26").Value = "3"
$ws.Range("BS26").Value = "1"
$ws.Range("BT26").Value = "1"
$ws.Range("BU26").Value = "2"
$ws.Range("BV26").Value = "2"
$ws.Range("BW26").Value = "2"

# Row 27
$ws.Range("A27").Value = "1"
$ws.Range("B27").Value = "1"
$ws.Range("C27").Value = "2"
$ws.Range("D27").Value = "1"
$ws.Range("E27").Value = "2"
$ws.Range("F27").Value = "2"
$ws.Range("G27").Value = "2"
$ws.Range("H27").Value = "2"
$ws.Range("I27").Value = "3"
$ws.Range("J27").Value = "2"
$ws.Range("K27").Value = "3"
$ws.Range("L27").Value = "3"
$ws.Range("M27").Value = "2"
$ws.Range("N27").Value = "3"
$ws.Range("O27").Value = "3"
$ws.Range("P27").Value = "1"
$ws.Range("Q27").Value = "1"
$ws.Range("R27").Value = "2"
$ws.Range("S27").Value = "1"
$ws.Range("T27").Value = "2"
$ws.Range("U27").Value = "3"
$ws.Range("V27").Value = "3"
$ws.Range("W27").Value = "3"
$ws.Range("X27").Value = "3"
$ws.Range("Y27").Value = "3"
$ws.Range("Z27").Value = "3"
$ws.Range("AA27").Value = "1"
$ws.Range("AB27").Value = "3"
$ws.Range("AC27").Value = "2"
$ws.Range("AD27").Value = "3"
$ws.Range("AE27").Value = "1"
$ws.Range("AF27").Value = "3"
$ws.Range("AG27").Value = "2"
$ws.Range("AH27").Value = "1"
$ws.Range("AI27").Value = "2"
$ws.Range("AJ27").Value = "2"
$ws.Range("AK27").Value = "2"
$ws.Range("AL27").Value = "3"
$ws.Range("AM27").Value = "1"
$ws.Range("AN27").Value = "2"
$ws.Range("AO27").Value = "2"
$ws.Range("AP27").Value = "3"
$ws.Range("AQ27").Value = "2"
$ws.Range("AR27").Value = "1"
$ws.Range("AS27").Value = "2"
$ws.Range("AT27").Value = "3"
$ws.Range("AU27").Value = "3"
$ws.Range("AV27").Value = "3"
$ws.Range("AW27").Value = "3"
$ws.Range("AX27").Value = "3"
$ws.Range("AY27").Value = "1"
$ws.Range("AZ27").Value = "3"
$ws.Range("BA27").Value = "2"
$ws.Range("BB27").Value = "1"
$ws.Range("BC27").Value = "3"
$ws.Range("BD27").Value = "2"
$ws.Range("BE27").Value = "2"
$ws.Range("BF27").Value = "2"
$ws.Range("BG27").Value = "2"
$ws.Range("BH27").Value = "2"
$ws.Range("BI27").Value = "3"
$ws.Range("BJ27").Value = "3"
$ws.Range("BK27").Value = "3"
$ws.Range("BL27").Value = "3"
$ws.Range("BM27").Value = "3"
$ws.Range("BN27").Value = "1"
$ws.Range("BO27").Value = "1"
$ws.Range("BP27").Value = "3"
$ws.Range("BQ27").Value = "1"
$ws.Range("BR27").Value = "3"
$ws.Range("BS27").Value = "1"
$ws.Range("BT27").Value = "1"
$ws.Range("BU27").Value = "2"
$ws.Range("BV27").Value = "1"
$ws.Range("BW27").Value = "2"

# Re-apply the same cell style/format used by the other response rows
# (copying it wipes out the temporary "@" text-format style created above).
$ws.Range("A24:BW24").Copy()
$ws.Range("A25:BW27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the selection left behind by the edit session.
$ws.Range("CB11").Select()
